# Apply the update described by the diff:
# A new data row is inserted at worksheet row 368 (pushing the existing
# rows 368-398 down to 369-399), containing a new price record for
# "Ají" variety "Inferno" dated 45013 (2023-03-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 368; this shifts all the
# existing rows (368..398) down by one (to 369..399) along with their
# values/styles, matching the diff.
$ws.Rows.Item(368).Insert()

# Populate the newly inserted row 368 with the new record values.
$ws.Cells.Item(368, 1).Value  = 8
$ws.Cells.Item(368, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(368, 3).Value  = "Coquimbo"
$ws.Cells.Item(368, 4).Value  = 45013
$ws.Cells.Item(368, 5).Value  = 4
$ws.Cells.Item(368, 6).Value  = 100112021
$ws.Cells.Item(368, 7).Value  = "Ají"
$ws.Cells.Item(368, 8).Value  = "Inferno"
$ws.Cells.Item(368, 9).Value  = "Primera"
$ws.Cells.Item(368, 10).Value = 480
$ws.Cells.Item(368, 11).Value = 11000
$ws.Cells.Item(368, 12).Value = 12000
$ws.Cells.Item(368, 13).Value = 11500
$ws.Cells.Item(368, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(368, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(368, 16).Value = 767
$ws.Cells.Item(368, 17).Value = 15
$ws.Cells.Item(368, 18).Value = "Hortaliza"

# Give the new row the same date style (s="2") as the rest of column D.
$ws.Cells.Item(368, 4).NumberFormat = $ws.Cells.Item(369, 4).NumberFormat
